# Fix for excelreader exceptions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "Sheet1" to "test"
$ws.Name = "test"

# Reset all columns to default width/style (Select All -> reset)
$ws.Cells.EntireColumn.AutoFit()
$ws.Columns.ColumnWidth = 11.53

# Center-align the header/data range B1:F3
$rng = $ws.Range("B1:F3")
$rng.HorizontalAlignment = -4108  # xlCenter
$rng.Font.Name = $rng.Font.Name

# Move selection
$ws.Range("G8").Select()
